$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1324.3
$ws.Range("I38").Value = 1249.2222
$ws.Range("J38").Value = 2000
$ws.Range("K38").Value = 3747.6666
$ws.Range("L38").Value = 6000
$ws.Range("M38").Value = -3375.6666
$ws.Range("N38").Value = -6744

$ws.Range("H111").Value = 1263.7646
$ws.Range("I111").Value = 775.7143
$ws.Range("J111").Value = 1605.4
$ws.Range("K111").Value = 2327.1429
$ws.Range("L111").Value = 4816.200000000001
$ws.Range("M111").Value = 739.8571000000002
$ws.Range("N111").Value = -10950.2

$ws.Range("H113").Value = 5168.5386
$ws.Range("I113").Value = 3393
$ws.Range("J113").Value = 6278.25
$ws.Range("K113").Value = 3393
$ws.Range("L113").Value = 6278.25
$ws.Range("M113").Value = -139

$ws.Range("H138").Value = 9377821
$ws.Range("I138").Value = 2085741.2
$ws.Range("J138").Value = 16669901
$ws.Range("K138").Value = 6257223.6
$ws.Range("L138").Value = 50009703
$ws.Range("M138").Value = -6252083.6
$ws.Range("N138").Value = -50019983

$ws.Range("H139").Value = 50780
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 50780
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 50780
$ws.Range("N139").Value = -61060

$ws.Range("H140").Value = 74675
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 74675
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 74675
$ws.Range("N140").Value = -85035

$ws.Range("H141").Value = 3920.1226
$ws.Range("I141").Value = 2198.2778
$ws.Range("J141").Value = 8688.308000000001
$ws.Range("K141").Value = 6594.8334
$ws.Range("L141").Value = 26064.924
$ws.Range("M141").Value = -1414.8334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 203606.6
$ws.Range("I2").Value = 203606.6
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 203606.6
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -203493.6

$ws.Range("H45").Value = 1112.8572
$ws.Range("I45").Value = 1098.3334
$ws.Range("J45").Value = 1200
$ws.Range("K45").Value = 1098.3334
$ws.Range("L45").Value = 1200
$ws.Range("M45").Value = -721.3334
$ws.Range("N45").Value = -1954

$ws.Range("H61").Value = 4066.5715
$ws.Range("I61").Value = 2946.1667
$ws.Range("J61").Value = 4906.875
$ws.Range("K61").Value = 2946.1667
$ws.Range("L61").Value = 4906.875
$ws.Range("M61").Value = -2734.1667

$ws.Range("H62").Value = 20000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 20000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 20000
$ws.Range("N62").Value = -21248

$ws.Range("H65").Value = 20000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 20000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 60000
$ws.Range("N65").Value = -66240

$ws.Range("H116").Value = 203606.6
$ws.Range("I116").Value = 203606.6
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 203606.6
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -201312.6

$ws.Range("H122").Value = 15537.5
$ws.Range("I122").Value = 19633.334
$ws.Range("J122").Value = 3250
$ws.Range("K122").Value = 58900.00199999999
$ws.Range("L122").Value = 9750
$ws.Range("M122").Value = -56450.00199999999

$ws.Range("H132").Value = 3078.162
$ws.Range("I132").Value = 2479.3928
$ws.Range("J132").Value = 4941
$ws.Range("K132").Value = 7438.178400000001
$ws.Range("L132").Value = 14823
$ws.Range("M132").Value = -4908.178400000001
$ws.Range("N132").Value = -19883

$ws.Range("H136").Value = 4066.5715
$ws.Range("I136").Value = 2946.1667
$ws.Range("J136").Value = 4906.875
$ws.Range("K136").Value = 8838.500100000001
$ws.Range("L136").Value = 14720.625
$ws.Range("M136").Value = -6288.500100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 203606.6
$ws.Range("I3").Value = 203606.6
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 203606.6
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -203492.6

$ws.Range("H94").Value = 2599.7
$ws.Range("I94").Value = 3286.6667
$ws.Range("J94").Value = 1569.25
$ws.Range("K94").Value = 3286.6667
$ws.Range("L94").Value = 1569.25
$ws.Range("M94").Value = -2835.6667

$ws.Range("H134").Value = 3116.4546
$ws.Range("I134").Value = 2335.2903
$ws.Range("J134").Value = 4979.231
$ws.Range("K134").Value = 7005.8709
$ws.Range("L134").Value = 14937.693
$ws.Range("M134").Value = -4470.8709
$ws.Range("N134").Value = -20007.693

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H114").Value = 30316.666
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 30316.666
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 30316.666
$ws.Range("N114").Value = -38994.666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 573.05
$ws.Range("I113").Value = 433.2857
$ws.Range("J113").Value = 648.3077
$ws.Range("K113").Value = 1299.8571
$ws.Range("L113").Value = 1944.9231
$ws.Range("M113").Value = 870.1428999999998
$ws.Range("N113").Value = -6284.9231

$ws.Range("H131").Value = 5377709
$ws.Range("I131").Value = 367.8
$ws.Range("J131").Value = 6411813
$ws.Range("K131").Value = 1103.4
$ws.Range("L131").Value = 19235439
$ws.Range("M131").Value = 3936.6
$ws.Range("N131").Value = -19245519

$ws.Range("H133").Value = 11971.429
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 11971.429
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 35914.287
$ws.Range("N133").Value = -46034.287
$ws.Range("M133").ClearContents()

$ws.Range("H136").Value = 2835.795
$ws.Range("I136").Value = 2586.6667
$ws.Range("J136").Value = 2856.5557
$ws.Range("K136").Value = 7760.000100000001
$ws.Range("L136").Value = 8569.667099999999
$ws.Range("M136").Value = -2660.000100000001
$ws.Range("N136").Value = -18769.6671

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 3000
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 3000
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 3000
$ws.Range("N9").Value = -3340
$ws.Range("M9").ClearContents()

$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2000
$ws.Range("N113").Value = -6340
$ws.Range("M113").ClearContents()

$ws.Range("H122").Value = 2766.6667
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2766.6667
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 8300.000100000001
$ws.Range("N122").Value = -13200.0001
$ws.Range("M122").ClearContents()

$ws.Range("H132").Value = 3713.3225
$ws.Range("I132").Value = 2758.3635
$ws.Range("J132").Value = 6047.6665
$ws.Range("K132").Value = 8275.0905
$ws.Range("L132").Value = 18142.9995
$ws.Range("M132").Value = -5745.0905
$ws.Range("N132").Value = -23202.9995

$ws.Range("H133").Value = 49682.855
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 49682.855
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 49682.855
$ws.Range("N133").Value = -59802.855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3188.5557
$ws.Range("I7").Value = 2710.1
$ws.Range("J7").Value = 3470
$ws.Range("K7").Value = 2710.1
$ws.Range("L7").Value = 3470
$ws.Range("M7").Value = -2598.1
$ws.Range("N7").Value = -3694

$ws.Range("H61").Value = 923.5294
$ws.Range("I61").Value = 616.2222
$ws.Range("J61").Value = 1269.25
$ws.Range("K61").Value = 616.2222
$ws.Range("L61").Value = 1269.25
$ws.Range("M61").Value = -414.2222
$ws.Range("N61").Value = -1673.25

$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H113").Value = 923.5294
$ws.Range("I113").Value = 616.2222
$ws.Range("J113").Value = 1269.25
$ws.Range("K113").Value = 616.2222
$ws.Range("L113").Value = 1269.25
$ws.Range("M113").Value = 1553.7778
$ws.Range("N113").Value = -5609.25

$ws.Range("H122").Value = 3804.7273
$ws.Range("I122").Value = 2701.3333
$ws.Range("J122").Value = 3978.9473
$ws.Range("K122").Value = 8103.999899999999
$ws.Range("L122").Value = 11936.8419
$ws.Range("M122").Value = -5653.999899999999

$ws.Range("H126").Value = 3188.5557
$ws.Range("I126").Value = 2710.1
$ws.Range("J126").Value = 3470
$ws.Range("K126").Value = 8130.299999999999
$ws.Range("L126").Value = 10410
$ws.Range("M126").Value = -5660.299999999999
$ws.Range("N126").Value = -15350

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1660.4
$ws.Range("I107").Value = 1660.4
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 4981.200000000001
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -3061.200000000001

$ws.Range("H122").Value = 1928.5
$ws.Range("I122").Value = 1527.8889
$ws.Range("J122").Value = 2443.5715
$ws.Range("K122").Value = 4583.6667
$ws.Range("L122").Value = 7330.7145
$ws.Range("M122").Value = -2133.6667
$ws.Range("N122").Value = -12230.7145

$ws.Range("H132").Value = 3024.932
$ws.Range("I132").Value = 2880.8484
$ws.Range("J132").Value = 3457.182
$ws.Range("K132").Value = 8642.5452
$ws.Range("L132").Value = 10371.546
$ws.Range("M132").Value = -6112.5452
$ws.Range("N132").Value = -15431.546
